# Generate Report for Handoff
# Update the "Latest Handoff" timestamps for the file
# 7895afa4-0fb9-461c-a484-9d353b3a4d53 (status: "Ready for handoff")
# to reflect a fresh handoff across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# Overview sheet: row 5 corresponds to 7895afa4-0fb9-461c-a484-9d353b3a4d53.md
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D5").Value = "2016-31-20 08:31:43"

# zh-cn sheet: row 5 is the same file's zh-cn handoff info
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E5").Value = "2016-03-20 08:31:40"

# de-de sheet: row 5 is the same file's de-de handoff info
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E5").Value = "2016-03-20 08:31:43"
